# Generate Report for Handback
#
# Marks the zh-cn and de-de localization rows as handed back: updates the
# "Status" text (shared across the Overview + per-locale sheets), fills in
# the "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" columns for each locale's two rows, and turns the new
# "Latest Handback File" entries into hyperlinks (same style/target as the
# existing "Latest Handoff File" links).

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

$mdTarget  = "721691ad-7b3d-415f-b730-8c9ee5d775ff.md"
$mdUrl     = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c97f16c6d1f80054b95482263bc19d0325b44b25/e2e/721691ad-7b3d-415f-b730-8c9ee5d775ff.md"
$hyperlinkColor = 15570276   # BGR long for RGB 6495ED (cornflower blue), matches the workbook's existing HyperLink style

# ---- Overview sheet: Status shown for both locales on both rows ----
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = $newStatus
$ov.Range("F2").Value = $newStatus
$ov.Range("E3").Value = $newStatus
$ov.Range("F3").Value = $newStatus

# ---- zh-cn sheet ----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = $newStatus
$zh.Range("C3").Value = $newStatus

$zh.Range("J2").Value = "721691ad-7b3d-415f-b730-8c9ee5d775ff.25efd87e78c42f047c0dc7ff0d780f1539a66c9d.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-25 09:07:47"
$zh.Range("J3").Value = "721691ad-7b3d-415f-b730-8c9ee5d775ff.25efd87e78c42f047c0dc7ff0d780f1539a66c9d.zh-cn.xlf"
$zh.Range("K3").Value = "2016-08-25 09:07:47"

$zhI2 = $zh.Range("I2")
$zh.Hyperlinks.Add($zhI2, $mdUrl, "", "", $mdTarget)
$zhI2.Font.Underline = $true
$zhI2.Font.Color = $hyperlinkColor

$zhI3 = $zh.Range("I3")
$zh.Hyperlinks.Add($zhI3, $mdUrl, "", "", $mdTarget)
$zhI3.Font.Underline = $true
$zhI3.Font.Color = $hyperlinkColor

# ---- de-de sheet ----
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = $newStatus
$de.Range("C3").Value = $newStatus

$de.Range("J2").Value = "721691ad-7b3d-415f-b730-8c9ee5d775ff.25efd87e78c42f047c0dc7ff0d780f1539a66c9d.de-de.xlf"
$de.Range("K2").Value = "2016-08-25 09:07:54"
$de.Range("J3").Value = "721691ad-7b3d-415f-b730-8c9ee5d775ff.25efd87e78c42f047c0dc7ff0d780f1539a66c9d.de-de.xlf"
$de.Range("K3").Value = "2016-08-25 09:07:54"

$deI2 = $de.Range("I2")
$de.Hyperlinks.Add($deI2, $mdUrl, "", "", $mdTarget)
$deI2.Font.Underline = $true
$deI2.Font.Color = $hyperlinkColor

$deI3 = $de.Range("I3")
$de.Hyperlinks.Add($deI3, $mdUrl, "", "", $mdTarget)
$deI3.Font.Underline = $true
$deI3.Font.Color = $hyperlinkColor

# ---- Column widths widened to fit the newly-populated target/handback file columns ----
$ov.Columns.Item(5).ColumnWidth = 29.9777047293527
$ov.Columns.Item(6).ColumnWidth = 29.9777047293527

$zh.Columns.Item(3).ColumnWidth = 29.9777047293527
$zh.Columns.Item(9).ColumnWidth = 40
$zh.Columns.Item(10).ColumnWidth = 40

$de.Columns.Item(3).ColumnWidth = 29.9777047293527
$de.Columns.Item(9).ColumnWidth = 40
$de.Columns.Item(10).ColumnWidth = 40

Write-Output "Handback report generated"
